# TemplateDefintion.id isn't a generated value anymore.
# TemplateDefinition.name is now TemplateDefinition.Id.
#
# On the "Configuration" sheet:
#   - Row 4 used to hold the "Name" variable; it now holds the "Id" variable
#     (value/description text updated accordingly).
#   - Row 8, which used to hold the auto-generated "Id" value, is removed
#     entirely since the id is no longer a generated value.
#   - The previously-active "Configuration" tab is no longer the active
#     sheet; "Variables" becomes the active sheet again (its default state).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

# Row 4: "Name" -> "Id", and drop the "...by this name or the id
# (defined below)." qualifier from the description now that the id is the
# only identifier.
$ws.Range("A4").Value = "Id"
$ws.Range("C4").Value = "Find a unique name shortly describing the functionality of this template, e. g. 'Employee contract'. You may refer this definition Excel file by this id."

# Remove the old, separately generated "Id" row (A8:C8 - "Id" /
# "JZpnpojeSuN5JDqtm9KZ" / "Please do not modify this value.").
$ws.Rows(8).Delete()

# Restore "Variables" as the active/selected sheet (it was "Configuration").
$wb.Worksheets.Item("Variables").Activate()
